# Auto-sync: insert a new quotation row (quotations sheet) and its
# associated item rows (items sheet), shifting subsequent rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "quotations": insert new row at position 8
# ---------------------------------------------------------------------
$wsQ = $wb.Worksheets.Item("quotations")
$wsQ.Rows.Item(8).Insert()

$wsQ.Cells.Item(8, 1).Value  = "NTUzNzQwNTMtYmRmMC00ZDAxLWI5ZjgtOGMyOTQzNjMyZDRmOjU3MDE2"   # id
$wsQ.Cells.Item(8, 2).Value  = "BJFPMUNTXG"                                                  # identifier
$wsQ.Cells.Item(8, 3).Value  = "MIX SOLUCOES AMBIENTAIS LTDA"                                # name
# description (D8) / additionalInformation (E8) stay blank
$wsQ.Cells.Item(8, 6).Value  = $false                                                        # archived
$wsQ.Cells.Item(8, 7).Value  = "'4984"                                                       # subtotal (text)
$wsQ.Cells.Item(8, 8).Value  = "'4984"                                                       # total (text)
$wsQ.Cells.Item(8, 9).Value  = "Pendente"                                                    # status
$wsQ.Cells.Item(8, 10).Value = "2025-10-15T21:09:42.233Z"                                    # expiredAt
# refusedAt (K8) stays blank
# approvedBy (L8) stays blank
$wsQ.Cells.Item(8, 13).Value = "Adriana Vieira Masini"                                       # createdBy
# approvedSignature (N8) stays blank
$wsQ.Cells.Item(8, 15).Value = "2025-10-08T21:12:33.362Z"                                    # createdAt
$wsQ.Cells.Item(8, 16).Value = "OTUxMWZiNzEtYjliOC00NTg4LWE5MTAtZmI2ZmQxZmZlZmNlOjU3MDE2"   # order.id
$wsQ.Cells.Item(8, 17).Value = "percentage"                                                  # discount.type
$wsQ.Cells.Item(8, 18).Value = "'0"                                                          # discount.value (text)
$wsQ.Cells.Item(8, 19).Value = "'0"                                                          # discount.total (text)
$wsQ.Cells.Item(8, 20).Value = "NDgzNDc2OTo1NzAxNg=="                                         # customer.id
$wsQ.Cells.Item(8, 21).Value = "pending"                                                     # status_original

# ---------------------------------------------------------------------
# Sheet "items": insert 3 new rows at position 38, linked to the new
# quotation created above
# ---------------------------------------------------------------------
$wsI = $wb.Worksheets.Item("items")
$wsI.Rows.Item(38).Insert()
$wsI.Rows.Item(38).Insert()
$wsI.Rows.Item(38).Insert()

# Row 38
$wsI.Cells.Item(38, 1).Value  = "NjA1ZWIxZDMtNWMxMC00NGE2LWEzZjEtYTZkZjM1MWVhNjQwOjU3MDE2"   # item_id
$wsI.Cells.Item(38, 2).Value  = 1                                                             # item_quantity
$wsI.Cells.Item(38, 3).Value  = 1793                                                          # item_total
# item_description (D38) stays blank
$wsI.Cells.Item(38, 5).Value  = 4                                                             # item_position
$wsI.Cells.Item(38, 6).Value  = "NTUzNzQwNTMtYmRmMC00ZDAxLWI5ZjgtOGMyOTQzNjMyZDRmOjU3MDE2"   # item_quotation_ref
$wsI.Cells.Item(38, 7).Value  = "NjdjMzI5NDAtMmU1Mi00MjQ1LTgxNGQtNjUyNWI3ZTQyNDU4OjU3MDE2"   # productService_id
$wsI.Cells.Item(38, 8).Value  = 1793                                                          # productService_value
$wsI.Cells.Item(38, 9).Value  = "product"                                                     # productService_type
$wsI.Cells.Item(38, 10).Value = "NTUzNzQwNTMtYmRmMC00ZDAxLWI5ZjgtOGMyOTQzNjMyZDRmOjU3MDE2"   # quotation_id

# Row 39
$wsI.Cells.Item(39, 1).Value  = "OTFkMjc2ZDEtOTkxOC00OGZlLWEyMWYtZGEwNDg3MDFiNzkxOjU3MDE2"   # item_id
$wsI.Cells.Item(39, 2).Value  = 1                                                             # item_quantity
$wsI.Cells.Item(39, 3).Value  = 350                                                           # item_total
$wsI.Cells.Item(39, 4).Value  = "Hora de trabalho ECO"                                        # item_description
$wsI.Cells.Item(39, 5).Value  = 4                                                             # item_position
$wsI.Cells.Item(39, 6).Value  = "NTUzNzQwNTMtYmRmMC00ZDAxLWI5ZjgtOGMyOTQzNjMyZDRmOjU3MDE2"   # item_quotation_ref
$wsI.Cells.Item(39, 7).Value  = "ODY3OTE5NTMtMDdjZi00YzM1LThkN2QtNDc5NzNmNzVkMGY0OjU3MDE2"   # productService_id
$wsI.Cells.Item(39, 8).Value  = 350                                                           # productService_value
$wsI.Cells.Item(39, 9).Value  = "service"                                                     # productService_type
$wsI.Cells.Item(39, 10).Value = "NTUzNzQwNTMtYmRmMC00ZDAxLWI5ZjgtOGMyOTQzNjMyZDRmOjU3MDE2"   # quotation_id

# Row 40
$wsI.Cells.Item(40, 1).Value  = "YTlhMGVhYzMtYjJhZS00OWUwLTg5YWQtNjdjYzIyMWUyZDZmOjU3MDE2"   # item_id
$wsI.Cells.Item(40, 2).Value  = 1                                                             # item_quantity
$wsI.Cells.Item(40, 3).Value  = 2841                                                          # item_total
# item_description (D40) stays blank
$wsI.Cells.Item(40, 5).Value  = 4                                                             # item_position
$wsI.Cells.Item(40, 6).Value  = "NTUzNzQwNTMtYmRmMC00ZDAxLWI5ZjgtOGMyOTQzNjMyZDRmOjU3MDE2"   # item_quotation_ref
$wsI.Cells.Item(40, 7).Value  = "MGQ3YTYzZmEtOGQyZS00YWNiLTljMWYtNTNiM2JkMzRmOTYwOjU3MDE2"   # productService_id
$wsI.Cells.Item(40, 8).Value  = 2841                                                          # productService_value
$wsI.Cells.Item(40, 9).Value  = "product"                                                     # productService_type
$wsI.Cells.Item(40, 10).Value = "NTUzNzQwNTMtYmRmMC00ZDAxLWI5ZjgtOGMyOTQzNjMyZDRmOjU3MDE2"   # quotation_id
